$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 1966.6666
$ws.Range("I52").Value = 450
$ws.Range("J52").Value = 5000
$ws.Range("K52").Value = 1350
$ws.Range("L52").Value = 15000
$ws.Range("M52").Value = -1190
$ws.Range("N52").Value = -15320
$ws.Range("H53").Value = 498
$ws.Range("I53").Value = 498
$ws.Range("K53").Value = 498
$ws.Range("M53").Value = 139
$ws.Range("H64").Value = 9500
$ws.Range("J64").Value = 9500
$ws.Range("L64").Value = 9500
$ws.Range("N64").Value = -9996
$ws.Range("H67").Value = 9500
$ws.Range("J67").Value = 9500
$ws.Range("L67").Value = 9500
$ws.Range("N67").Value = -11216
$ws.Range("H68").Value = 19000
$ws.Range("I68").Value = 19000
$ws.Range("K68").Value = 19000
$ws.Range("M68").Value = -18251
$ws.Range("H71").Value = 19000
$ws.Range("I71").Value = 19000
$ws.Range("K71").Value = 57000
$ws.Range("M71").Value = -53256
$ws.Range("H86").Value = 900
$ws.Range("J86").Value = 900
$ws.Range("L86").Value = 900
$ws.Range("N86").Value = -3146
$ws.Range("H89").Value = 900
$ws.Range("J89").Value = 900
$ws.Range("L89").Value = 4500
$ws.Range("N89").Value = -15732
$ws.Range("H92").Value = 1113.2858
$ws.Range("I92").Value = 449.5
$ws.Range("K92").Value = 449.5
$ws.Range("M92").Value = 798.5
$ws.Range("H96").Value = 2183.3333
$ws.Range("I96").Value = 4450
$ws.Range("K96").Value = 13350
$ws.Range("M96").Value = -11977
$ws.Range("H125").Value = 32
$ws.Range("I125").Value = 32
$ws.Range("K125").Value = 288
$ws.Range("M125").Value = 2172
$ws.Range("H137").Value = 4993
$ws.Range("J137").Value = 4993
$ws.Range("L137").Value = 14979
$ws.Range("N137").Value = -20079

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1655.3636
$ws.Range("I32").Value = 1655.3636
$ws.Range("K32").Value = 1655.3636
$ws.Range("M32").Value = -1368.3636
$ws.Range("H76").Value = 19562.666
$ws.Range("J76").Value = 19562.666
$ws.Range("L76").Value = 19562.666
$ws.Range("N76").Value = -20238.666
$ws.Range("H79").Value = 19562.666
$ws.Range("J79").Value = 19562.666
$ws.Range("L79").Value = 19562.666
$ws.Range("N79").Value = -21902.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2556.6365
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5792.857
$ws.Range("I31").Value = 1991
$ws.Range("J31").Value = 6426.5
$ws.Range("K31").Value = 1991
$ws.Range("L31").Value = 6426.5
$ws.Range("M31").Value = -1696
$ws.Range("N31").Value = -7016.5
$ws.Range("H34").Value = 5792.857
$ws.Range("I34").Value = 1991
$ws.Range("J34").Value = 6426.5
$ws.Range("K34").Value = 1991
$ws.Range("L34").Value = 6426.5
$ws.Range("M34").Value = -1789
$ws.Range("N34").Value = -6830.5
$ws.Range("H107").Value = 567.5
$ws.Range("I107").Value = 481
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 481
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1439
$ws.Range("N107").Value = -4840
$ws.Range("H134").Value = 2258.6667
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 2258.6667
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 6776.000100000001
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -11846.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 54.6
$ws.Range("I2").Value = 19.125
$ws.Range("J2").Value = 78.25
$ws.Range("K2").Value = 114.75
$ws.Range("L2").Value = 469.5
$ws.Range("M2").Value = -1.75
$ws.Range("N2").Value = -695.5
$ws.Range("H38").Value = 48.6
$ws.Range("I38").Value = 26.4
$ws.Range("J38").Value = 70.8
$ws.Range("K38").Value = 79.19999999999999
$ws.Range("L38").Value = 212.4
$ws.Range("N38").Value = -906.4
$ws.Range("M38").Value = 267.8
$ws.Range("H70").Value = 6400
$ws.Range("I70").Value = 6400
$ws.Range("K70").Value = 19200
$ws.Range("M70").Value = -18885
$ws.Range("H73").Value = 6400
$ws.Range("I73").Value = 6400
$ws.Range("K73").Value = 19200
$ws.Range("M73").Value = -18108

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 21000
$ws.Range("J26").Value = 21000
$ws.Range("L26").Value = 21000
$ws.Range("N26").Value = -21560
$ws.Range("H50").Value = 21000
$ws.Range("J50").Value = 21000
$ws.Range("L50").Value = 21000
$ws.Range("N50").Value = -21996
$ws.Range("H116").Value = 99995
$ws.Range("J116").Value = 99995
$ws.Range("L116").Value = 99995
$ws.Range("N116").Value = -109173

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9665
$ws.Range("I7").Value = 9497.5
$ws.Range("K7").Value = 9497.5
$ws.Range("M7").Value = -9385.5
$ws.Range("H40").Value = 3647.1667
$ws.Range("I40").Value = 3376.7
$ws.Range("K40").Value = 3376.7
$ws.Range("M40").Value = -3240.7
$ws.Range("H68").Value = 2599.9
$ws.Range("I68").Value = 2599.9
$ws.Range("K68").Value = 2599.9
$ws.Range("M68").Value = -1850.9
$ws.Range("H71").Value = 2599.9
$ws.Range("I71").Value = 2599.9
$ws.Range("K71").Value = 12999.5
$ws.Range("M71").Value = -9255.5
$ws.Range("H82").Value = 2368.625
$ws.Range("I82").Value = 2158.3333
$ws.Range("J82").Value = 2999.5
$ws.Range("K82").Value = 2158.3333
$ws.Range("L82").Value = 2999.5
$ws.Range("M82").Value = -1797.3333
$ws.Range("N82").Value = -3721.5
$ws.Range("H85").Value = 2368.625
$ws.Range("I85").Value = 2158.3333
$ws.Range("J85").Value = 2999.5
$ws.Range("K85").Value = 2158.3333
$ws.Range("L85").Value = 2999.5
$ws.Range("M85").Value = -910.3332999999998
$ws.Range("N85").Value = -5495.5
$ws.Range("H126").Value = 9665
$ws.Range("I126").Value = 9497.5
$ws.Range("K126").Value = 28492.5
$ws.Range("M126").Value = -26022.5
$ws.Range("H132").Value = 12500
$ws.Range("I132").Value = 10000
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 30000
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -27470
$ws.Range("N132").Value = -50060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 90162.5
$ws.Range("J86").Value = 90162.5
$ws.Range("L86").Value = 90162.5
$ws.Range("N86").Value = -92408.5
$ws.Range("H89").Value = 90162.5
$ws.Range("J89").Value = 90162.5
$ws.Range("L89").Value = 450812.5
$ws.Range("N89").Value = -462044.5
$ws.Range("H132").Value = 2889.111
$ws.Range("I132").Value = 3250.5
$ws.Range("J132").Value = 2600
$ws.Range("K132").Value = 9751.5
$ws.Range("L132").Value = 7800
$ws.Range("M132").Value = -7221.5
$ws.Range("N132").Value = -12860
